{"js": "// Replace the division-problem text runs in the table with new values.\n// Each \"old\" string is a unique, literal substring (a whole w:t run's\n// text), so a plain body-wide search-and-replace is safe and unambiguous.\nconst replacements = [\n  [\"517\u00f76=\", \"538\u00f79=\"],\n  [\"766\u00f72=\", \"461\u00f73=\"],\n  [\"930\u00f76=\", \"414\u00f75=\"],\n  [\"631\u00f77=\", \"946\u00f77=\"],\n  [\"623\u00f73=\", \"245\u00f75=\"],\n  [\"217\u00f73=\", \"320\u00f72=\"],\n  [\"701\u00f79=\", \"245\u00f78=\"],\n  [\"825\u00f77=\", \"781\u00f77=\"],\n  [\"277\u00f74=\", \"983\u00f78=\"],\n  [\"475\u00f75=\", \"886\u00f76=\"],\n  [\"462\u00f73=\", \"202\u00f77=\"],\n  [\"741\u00f74=\", \"155\u00f73=\"],\n  [\"318\u00f79=\", \"192\u00f72=\"],\n  [\"325\u00f75=\", \"103\u00f74=\"],\n  [\"951\u00f74=\", \"462\u00f78=\"],\n  [\"204\u00f76=\", \"539\u00f75=\"],\n  [\"443\u00f73=\", \"698\u00f79=\"],\n  [\"152\u00f72=\", \"354\u00f76=\"],\n  [\"247\u00f76=\", \"989\u00f73=\"],\n  [\"176\u00f72=\", \"228\u00f79=\"],\n  [\"984\u00f72=\", \"419\u00f77=\"],\n  [\"618\u00f77=\", \"735\u00f77=\"],\n  [\"697\u00f78=\", \"107\u00f75=\"],\n  [\"286\u00f78=\", \"404\u00f76=\"],\n  [\"621\u00f75=\", \"267\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text runs in the table with new values.\n# Each \"old\" string is a unique, literal whole-run text, so a simple\n# Find/Replace (wdReplaceAll) across the document body is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"517\u00f76=\"; New = \"538\u00f79=\"},\n    @{Old = \"766\u00f72=\"; New = \"461\u00f73=\"},\n    @{Old = \"930\u00f76=\"; New = \"414\u00f75=\"},\n    @{Old = \"631\u00f77=\"; New = \"946\u00f77=\"},\n    @{Old = \"623\u00f73=\"; New = \"245\u00f75=\"},\n    @{Old = \"217\u00f73=\"; New = \"320\u00f72=\"},\n    @{Old = \"701\u00f79=\"; New = \"245\u00f78=\"},\n    @{Old = \"825\u00f77=\"; New = \"781\u00f77=\"},\n    @{Old = \"277\u00f74=\"; New = \"983\u00f78=\"},\n    @{Old = \"475\u00f75=\"; New = \"886\u00f76=\"},\n    @{Old = \"462\u00f73=\"; New = \"202\u00f77=\"},\n    @{Old = \"741\u00f74=\"; New = \"155\u00f73=\"},\n    @{Old = \"318\u00f79=\"; New = \"192\u00f72=\"},\n    @{Old = \"325\u00f75=\"; New = \"103\u00f74=\"},\n    @{Old = \"951\u00f74=\"; New = \"462\u00f78=\"},\n    @{Old = \"204\u00f76=\"; New = \"539\u00f75=\"},\n    @{Old = \"443\u00f73=\"; New = \"698\u00f79=\"},\n    @{Old = \"152\u00f72=\"; New = \"354\u00f76=\"},\n    @{Old = \"247\u00f76=\"; New = \"989\u00f73=\"},\n    @{Old = \"176\u00f72=\"; New = \"228\u00f79=\"},\n    @{Old = \"984\u00f72=\"; New = \"419\u00f77=\"},\n    @{Old = \"618\u00f77=\"; New = \"735\u00f77=\"},\n    @{Old = \"697\u00f78=\"; New = \"107\u00f75=\"},\n    @{Old = \"286\u00f78=\"; New = \"404\u00f76=\"},\n    @{Old = \"621\u00f75=\"; New = \"267\u00f73=\"}\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $rep.New\n    $find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n}\n\n$d.Save()\n"}
